$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test2")

# Row 4: the sheet stores every cell as text (even numeric-looking ones,
# e.g. "7", "39", "15" -> now "21", "25", "45"). A leading apostrophe keeps
# Excel from re-interpreting the digits as a Number when we type them in.
$ws.Range("B4").Value = "'21"
$ws.Range("C4").Value = "'25"
$ws.Range("D4").Value = "'45"

# The last sample row (row 6: 7 / 36 / 43 under the HOUR/MINUTE/SECOND
# header in row 5) is removed entirely.
$ws.Rows.Item(6).Delete()
